$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; all rows 33..109 shift down to 34..110.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44526
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112045
$ws.Range("G33").Value = "Zapallo"
$ws.Range("H33").Value = "Paine"
$ws.Range("I33").Value = "1a (guarda)"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 220
$ws.Range("L33").Value = 250
$ws.Range("M33").Value = 235
$ws.Range("N33").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O33").Value = "Región de O'Higgins"
$ws.Range("P33").Value = 235
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"
